# "6 hours by turn fix"
# The weekly schedule grid gets re-balanced: a few existing slots swap their
# class label, and a whole new 15:50/17:30/18:20 tail of slots is appended
# (pushing the old 16:40 row down one position) so the day now covers a
# full 6-hour-by-turn span.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the subject that had drifted to the wrong period -------------
$ws.Range("D3").Value = "-"
$ws.Range("D6").Value = "MEC-2B-Resistencia mecanica"
$ws.Range("F6").Value = "-"
$ws.Range("F7").Value = "MEC-2B-Resistencia mecanica"

# --- Lunch break actually starts at 12:20, not 13:00 -------------------
$ws.Range("B8:F8").Value = "-"
$ws.Range("A9").Value = "12:20"
$ws.Range("B9:F9").Value = "Almoço"

# --- Shift the remaining afternoon period start-times back one slot ---
$ws.Range("A10").Value = "13:00"
$ws.Range("A11").Value = "13:50"
$ws.Range("A12").Value = "14:40"
$ws.Range("B12:F12").Value = "-"
$ws.Range("A13").Value = "15:30"
$ws.Range("B13:F13").Value = "Intervalo"

# --- Insert the new 15:50 slot, pushing the existing 16:40 row to 15 --
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "15:50"
$ws.Range("B14:F14").Value = "-"

# --- Append the new end-of-day slots: 17:30 and 18:20 -------------------
$ws.Range("A16").Value = "17:30"
$ws.Range("B16:F16").Value = "-"
$ws.Range("A17").Value = "18:20"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
